# Nasar - 09Nov changes
#
# Updates the "iOS" locator sheet:
#   - Row 30 is repurposed from the "Your accounts" tab locator to the
#     new "show all" down-arrow image locator.
#   - Three brand-new locator rows are appended (37, 38, 39 in 1-based
#     Excel terms -> rows 36-38 below the previous last row, plus one
#     trailing blank-ish row that only carries the VALUE_TYPE column).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("iOS")
$ws.Activate()

# --- Row 30: replace existing key/value pair -----------------------------
# ELEMENT_KEY first, then ELEMENT_VALUE (matches how the row was edited).
$ws.Range("A30").Value = "PORTFOLIO_SHOWALL_DOWN_ARROW"
$ws.Range("B30").Value = "ic_expand_pill"

# --- New row 36: portfolio performance arrow ------------------------------
# ELEMENT_VALUE is entered before ELEMENT_KEY for the newly appended rows.
$ws.Range("B36").Value = "portfolio-performance-increase-arrow"
$ws.Range("A36").Value = "PORTFOLIO_PERFORMANCE_ARROW"
$ws.Range("E36").Value = "device-accessibilityid"

# --- New row 37: portfolio summary value change text ----------------------
$ws.Range("B37").Value = "portfolio-summary-value-change"
$ws.Range("A37").Value = "PORTFOLIO_SUMMARY_VALUE_CHANGE_TXT"
$ws.Range("E37").Value = "device-accessibilityid"

# --- New row 38: portfolio summary "refreshed at" text ---------------------
$ws.Range("B38").Value = "portfolio-summary-refreshed-at-text"
$ws.Range("A38").Value = "PORTFOLIO_SUMMARY_REFRESH_TIME_TXT"
$ws.Range("E38").Value = "device-accessibilityid"

# --- New row 39: VALUE_TYPE only -------------------------------------------
$ws.Range("E39").Value = "device-accessibilityid"

# --- View state: selection moves to B13, viewport scrolls up to row 4 -----
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B13").Select()
